# Adding CRM test from 10/5/2019 for backcalculation
# Appends two new data rows (54 and 55) to Sheet1, mirroring the existing
# row layout/format, and updates the active selection to E46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats, styles) of the last existing data
# row (53) down into the two new rows so the new date cells pick up the
# same date style etc. without inventing new number formats.
$ws.Range("A53:F53").Copy()
$ws.Range("A54:F54").PasteSpecial(-4122)
$ws.Range("A55:F55").PasteSpecial(-4122)

# --- Row 54: opend crm (10/5/2019 ---
$ws.Range("A54").Value = 43743
$ws.Range("B54").Value = 2286.5546243738299
$ws.Range("C54").Value = 2207.0300000000002
$ws.Range("E54").Value = 169
$ws.Range("F54").Value = "opend crm (10/5/2019"

# --- Row 55: opened crm (10/5/2019 ---
$ws.Range("A55").Value = 43743
$ws.Range("B55").Value = 2281.5549225887498
$ws.Range("C55").Value = 2207.0300000000002
$ws.Range("E55").Value = 169
$ws.Range("F55").Value = "opened crm (10/5/2019"

# Fill the % off formula down both new rows as a single shared formula
# group, matching the pattern used throughout the sheet.
$ws.Range("D54:D55").Formula = "=100*(B54-C54)/C54"

# Move the selection to match the post-edit cursor position.
$ws.Range("E46").Select() | Out-Null
